# Horarios Línea 141 - actualización 05:31:23
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with the
# latest scrape: refreshes the "Última actualización" / "Total filas" header
# cells, inserts a newly-seen arrival (215B_EL PATO @ 05:35) ahead of the
# existing rows, and appends the batch of newly scraped rows at the bottom of
# each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header refresh
$ws1.Cells.Item(2,1).Value = "Última actualización: 05:31:23"
$ws1.Cells.Item(3,1).Value = "Total filas: 29"

# A new row (215B_EL PATO, 4 min) was scraped ahead of the previous "15_ABASTO"
# row, so insert a blank row at 14 and push everything below it down by one.
$ws1.Rows.Item(14).Insert()

$ws1.Cells.Item(14,1).Value = "05:31:23"
$ws1.Cells.Item(14,2).Value = "05:35"
$ws1.Cells.Item(14,3).Value = "215B_EL PATO"
$ws1.Cells.Item(14,4).Value = 4
$ws1.Cells.Item(14,5).Value = "LP1912"

# Newly scraped rows appended at the bottom (rows 27-34)
$ws1.Cells.Item(27,1).Value = "05:31:23"
$ws1.Cells.Item(27,2).Value = "06:59"
$ws1.Cells.Item(27,3).Value = "14_ABASTO"
$ws1.Cells.Item(27,4).Value = 88
$ws1.Cells.Item(27,5).Value = "LP1912"

$ws1.Cells.Item(28,1).Value = "05:31:23"
$ws1.Cells.Item(28,2).Value = "07:04"
$ws1.Cells.Item(28,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(28,4).Value = 93
$ws1.Cells.Item(28,5).Value = "LP1912"

$ws1.Cells.Item(29,1).Value = "05:31:23"
$ws1.Cells.Item(29,2).Value = "07:05"
$ws1.Cells.Item(29,3).Value = "15_ABASTO"
$ws1.Cells.Item(29,4).Value = 94
$ws1.Cells.Item(29,5).Value = "LP1912"

$ws1.Cells.Item(30,1).Value = "05:31:23"
$ws1.Cells.Item(30,2).Value = "07:07"
$ws1.Cells.Item(30,3).Value = "225_GOMEZ"
$ws1.Cells.Item(30,4).Value = 96
$ws1.Cells.Item(30,5).Value = "LP1912"

$ws1.Cells.Item(31,1).Value = "05:31:23"
$ws1.Cells.Item(31,2).Value = "07:11"
$ws1.Cells.Item(31,3).Value = "215A_EL PATO"
$ws1.Cells.Item(31,4).Value = 100
$ws1.Cells.Item(31,5).Value = "LP1912"

$ws1.Cells.Item(32,1).Value = "05:31:23"
$ws1.Cells.Item(32,2).Value = "07:15"
$ws1.Cells.Item(32,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(32,4).Value = 104
$ws1.Cells.Item(32,5).Value = "LP1912"

$ws1.Cells.Item(33,1).Value = "05:31:23"
$ws1.Cells.Item(33,2).Value = "07:21"
$ws1.Cells.Item(33,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(33,4).Value = 110
$ws1.Cells.Item(33,5).Value = "LP1912"

$ws1.Cells.Item(34,1).Value = "05:31:23"
$ws1.Cells.Item(34,2).Value = "07:23"
$ws1.Cells.Item(34,3).Value = "10_OLMOS"
$ws1.Cells.Item(34,4).Value = 112
$ws1.Cells.Item(34,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

# Header refresh
$ws2.Cells.Item(2,1).Value = "Última actualización: 05:31:23"
$ws2.Cells.Item(3,1).Value = "Total filas: 7"

# The previous row 9 (215A_EL PATO @ 06:11) is replaced by the freshly scraped
# 215B_EL PATO @ 05:35 row; its old content moves down into a newly inserted
# row 10.
$ws2.Rows.Item(10).Insert()

$ws2.Cells.Item(9,1).Value = "05:31:23"
$ws2.Cells.Item(9,2).Value = "05:35"
$ws2.Cells.Item(9,3).Value = "215B_EL PATO"
$ws2.Cells.Item(9,4).Value = 4
$ws2.Cells.Item(9,5).Value = "LP1912"

$ws2.Cells.Item(10,1).Value = "04:28:33"
$ws2.Cells.Item(10,2).Value = "06:11"
$ws2.Cells.Item(10,3).Value = "215A_EL PATO"
$ws2.Cells.Item(10,4).Value = 103
$ws2.Cells.Item(10,5).Value = "LP1912"

# Newly scraped row appended at the bottom
$ws2.Cells.Item(12,1).Value = "05:31:23"
$ws2.Cells.Item(12,2).Value = "07:11"
$ws2.Cells.Item(12,3).Value = "215A_EL PATO"
$ws2.Cells.Item(12,4).Value = 100
$ws2.Cells.Item(12,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

# Header refresh
$ws3.Cells.Item(2,1).Value = "Última actualización: 05:31:23"
$ws3.Cells.Item(3,1).Value = "Total filas: 7"

# Newly scraped row appended at the bottom
$ws3.Cells.Item(12,1).Value = "05:31:23"
$ws3.Cells.Item(12,2).Value = "07:00"
$ws3.Cells.Item(12,3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(12,4).Value = 89
$ws3.Cells.Item(12,5).Value = "L6173"
